# Adds %DiffH, %DiffD, %DiffA columns (AF, AG, AH) computed as
# (DiffX / YtrueX) * 100 for every data row, mirroring the existing
# DiffH/DiffD/DiffA columns (AC/AD/AE) next to YtrueH/YtrueD/YtrueA (Z/AA/AB).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: Z=26 (YtrueH), AA=27 (YtrueD), AB=28 (YtrueA)
#                 AC=29 (DiffH),  AD=30 (DiffD),  AE=31 (DiffA)
#                 AF=32 (%DiffH), AG=33 (%DiffD), AH=34 (%DiffA)

$lastRow = 150

# --- Header row (row 1): copy the formatting of the DiffH/DiffD/DiffA
#     header cells so the new headers match the existing bold/bordered style.
$ws.Cells.Item(1, 29).Copy()
$ws.Cells.Item(1, 32).PasteSpecial(-4122)
$ws.Cells.Item(1, 32).Value2 = "%DiffH"

$ws.Cells.Item(1, 30).Copy()
$ws.Cells.Item(1, 33).PasteSpecial(-4122)
$ws.Cells.Item(1, 33).Value2 = "%DiffD"

$ws.Cells.Item(1, 31).Copy()
$ws.Cells.Item(1, 34).PasteSpecial(-4122)
$ws.Cells.Item(1, 34).Value2 = "%DiffA"

# --- Data rows: compute %Diff = Diff / Ytrue * 100 for each of H/D/A.
for ($r = 2; $r -le $lastRow; $r++) {
    $YtrueH = $ws.Cells.Item($r, 26).Value2
    $YtrueD = $ws.Cells.Item($r, 27).Value2
    $YtrueA = $ws.Cells.Item($r, 28).Value2

    $DiffH = $ws.Cells.Item($r, 29).Value2
    $DiffD = $ws.Cells.Item($r, 30).Value2
    $DiffA = $ws.Cells.Item($r, 31).Value2

    $ws.Cells.Item($r, 32).Value2 = ($DiffH / $YtrueH) * 100
    $ws.Cells.Item($r, 33).Value2 = ($DiffD / $YtrueD) * 100
    $ws.Cells.Item($r, 34).Value2 = ($DiffA / $YtrueA) * 100
}

Write-Host "Added %DiffH/%DiffD/%DiffA columns (AF:AH) for rows 2-$lastRow"
